# Raw and Clean Data from SSA for August 4th
# Appends one new day (2020-08-04, Excel serial 44047) of data to the
# daily tracking workbook: out_vars, dates_dx, dates_sx, dates_deaths all
# get a new row 66; control_obs gets a new BN column (one more day); and
# a couple of cosmetic selection / active-cell updates follow along.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# out_vars: append row 66 with the day's raw headline numbers
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("out_vars")
$ws1.Activate()

$ws1.Range("A65:J65").Copy()
$ws1.Range("A66:J66").PasteSpecial(-4122)   # xlPasteFormats

$ws1.Range("A66").Value = 44047
$ws1.Range("B66").Value = 449961
$ws1.Range("C66").Value = 493873
$ws1.Range("D66").Value = 82460
$ws1.Range("E66").Value = 48869
$ws1.Range("F66").Value = 26.864550483264104
$ws1.Range("G66").Value = 120880
$ws1.Range("H66").Value = 9669
$ws1.Range("I66").Value = 11605
$ws1.Range("J66").Value = 1026294

$ws1.Range("A66").Select()

# ---------------------------------------------------------------
# dates_dx: append row 66, extend the A61:A66 highlighted-block style
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("dates_dx")
$ws2.Activate()

$ws2.Range("A65:L65").Copy()
$ws2.Range("A66:L66").PasteSpecial(-4122)

$ws2.Range("A60").Copy()
$ws2.Range("A61:A66").PasteSpecial(-4122)

$ws2.Range("A66").Value = 44047
$ws2.Range("B66").Value = 0
$ws2.Range("C66").Value = 1
$ws2.Range("D66").Value = 0
$ws2.Range("E66").Value = 0
$ws2.Range("F66").Value = 1
$ws2.Range("G66").Value = 0
$ws2.Range("H66").Value = 0
$ws2.Range("I66").Value = 0
$ws2.Range("J66").Value = 0
$ws2.Range("K66").Value = 0
$ws2.Range("L66").Value = 4

$ws2.Range("A60:A66").Select()

# ---------------------------------------------------------------
# dates_sx: append row 66, extend the A61:A66 highlighted-block style
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("dates_sx")
$ws3.Activate()

$ws3.Range("A65:N65").Copy()
$ws3.Range("A66:N66").PasteSpecial(-4122)

$ws3.Range("A60").Copy()
$ws3.Range("A61:A66").PasteSpecial(-4122)

$ws3.Range("A66").Value = 44047
$ws3.Range("B66").Value = 0
$ws3.Range("C66").Value = 1
$ws3.Range("D66").Value = 0
$ws3.Range("E66").Value = 0
$ws3.Range("F66").Value = 0
$ws3.Range("G66").Value = 0
$ws3.Range("H66").Value = 1
$ws3.Range("I66").Value = 0
$ws3.Range("J66").Value = 0
$ws3.Range("K66").Value = 1
$ws3.Range("L66").Value = 0
$ws3.Range("M66").Value = 0
$ws3.Range("N66").Value = 0

$ws3.Range("C66").Select()

# ---------------------------------------------------------------
# dates_deaths: append row 66, extend the A61:A66 highlighted-block style
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("dates_deaths")
$ws4.Activate()

$ws4.Range("A65:J65").Copy()
$ws4.Range("A66:J66").PasteSpecial(-4122)

$ws4.Range("A60").Copy()
$ws4.Range("A61:A66").PasteSpecial(-4122)

$ws4.Range("A66").Value = 44047
$ws4.Range("B66").Value = 0
$ws4.Range("C66").Value = 0
$ws4.Range("D66").Value = 0
$ws4.Range("E66").Value = 0
$ws4.Range("F66").Value = 2
$ws4.Range("G66").Value = 1
$ws4.Range("H66").Value = 1
$ws4.Range("I66").Value = 1
$ws4.Range("J66").Value = 2

$ws4.Range("K66").Select()

# ---------------------------------------------------------------
# control_obs: new BN column (one more day, 2020-08-04) of the big
# running tally, plus the BM15 correction that the new day's data
# revealed.
# ---------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("control_obs")
$ws5.Activate()

$ws5.Range("BN1").Value = 44047
$ws5.Range("BN2").Value = 4994
$ws5.Range("BN3").Value = 4793
$ws5.Range("BN4").Value = 4793
$ws5.Range("BN5").Value = 4793
$ws5.Range("BN6").Value = 4793
$ws5.Range("BN7").Value = 4103
$ws5.Range("BN8").Value = 6734
$ws5.Range("BN10").Value = 205
$ws5.Range("BN11").Value = 205
$ws5.Range("BN12").Value = 205
$ws5.Range("BN13").Value = 205
$ws5.Range("BN14").Value = 205

# row 15 correction: BM15 was mis-entered, fix it and populate BN15
$ws5.Range("BM15").Value = 139
$ws5.Range("BN15").Value = 140

$ws5.Range("BN16").Value = 217

# BN18 has no explicit style in the source (matches BF18/BG18, also unstyled)
$ws5.Range("BN18").Value = 1164

# extend the running SUM() total one more column
$ws5.Range("BM20").Copy()
$ws5.Range("BN20").PasteSpecial(-4122)
$ws5.Range("BN20").Formula = "=SUM(BN2:BN18)"

$ws5.Range("BO23").Select()

# restore control_obs as the active/selected tab (it was already the
# active sheet before this script ran)
$ws5.Activate()

Write-Output "edit applied"
